$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.181.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.03'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.23'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4632'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3741'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07384'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8664'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.56'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.814.15'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.651'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.378'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07084'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.63'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008736'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.87'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.184.88'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.303'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.30%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.049.96'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.929'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.68'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.217'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.52'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.261'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.77'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08882'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7721'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.170'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.508'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.927'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.63%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01958'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05235'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.240'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.907'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.365'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +20.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5268'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1679'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.569'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5016'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.35'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.23'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.666'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.27%  '
